# Updated symbol list on Tue Jan 31 12:00:02 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) values for the crypto
# symbol rows on the active sheet to the latest scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.87"
$ws.Range("E2").Value = "'0.66%"
$ws.Range("D3").Value = "'37.35"
$ws.Range("E3").Value = "'-0.53%"
$ws.Range("D4").Value = "'5.129"
$ws.Range("E4").Value = "'1.53%"
$ws.Range("D5").Value = "'0.07778"
$ws.Range("E5").Value = "'-1.45%"
$ws.Range("D6").Value = "'8.190"
$ws.Range("E6").Value = "'-0.37%"
$ws.Range("D7").Value = "'1.874"
$ws.Range("E7").Value = "'-7.30%"
$ws.Range("D8").Value = "'2.882"
$ws.Range("E8").Value = "'-7.76%"
$ws.Range("D9").Value = "'0.9196"
$ws.Range("E10").Value = "'-6.83%"
$ws.Range("D11").Value = "'0.1895"
$ws.Range("E11").Value = "'0.11%"
$ws.Range("D12").Value = "'0.09387"
$ws.Range("E12").Value = "'7.83%"
$ws.Range("D13").Value = "'0.03409"
$ws.Range("E13").Value = "'-0.74%"
$ws.Range("D14").Value = "'0.09690"
$ws.Range("E14").Value = "'-0.49%"
$ws.Range("D15").Value = "'0.001374"
$ws.Range("E15").Value = "'-1.30%"
$ws.Range("D16").Value = "'0.005775"
$ws.Range("E16").Value = "'-3.64%"
$ws.Range("D17").Value = "'3.534"
$ws.Range("E17").Value = "'-1.11%"
$ws.Range("D18").Value = "'4.405"
$ws.Range("E18").Value = "'1.02%"
$ws.Range("D19").Value = "'0.3401"
$ws.Range("E19").Value = "'-1.02%"
$ws.Range("D20").Value = "'5.255"
$ws.Range("E20").Value = "'4.61%"
$ws.Range("D21").Value = "'0.1267"
$ws.Range("E21").Value = "'-1.45%"
$ws.Range("D22").Value = "'0.2591"
$ws.Range("E22").Value = "'2.73%"
$ws.Range("E23").Value = "'180.26%"
$ws.Range("D24").Value = "'0.04325"
$ws.Range("E24").Value = "'-0.42%"
$ws.Range("E25").Value = "'-2.10%"
$ws.Range("D26").Value = "'0.004254"
$ws.Range("E26").Value = "'-7.87%"
$ws.Range("D27").Value = "'0.0001301"
$ws.Range("E27").Value = "'-63.80%"
$ws.Range("D39").Value = "'0.02072"
$ws.Range("E39").Value = "'-8.25%"
$ws.Range("D40").Value = "'0.05034"
$ws.Range("E40").Value = "'0.36%"
$ws.Range("D41").Value = "'0.007664"
$ws.Range("E41").Value = "'1.98%"
$ws.Range("D42").Value = "'0.009816"
$ws.Range("E42").Value = "'-1.08%"
$ws.Range("D43").Value = "'0.1344"
$ws.Range("E43").Value = "'-1.11%"
$ws.Range("D44").Value = "'0.002171"
$ws.Range("E44").Value = "'7.02%"
$ws.Range("D45").Value = "'0.008735"
$ws.Range("E45").Value = "'2.13%"
$ws.Range("D46").Value = "'0.00006707"
$ws.Range("E46").Value = "'4.63%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.39%"
$ws.Range("D48").Value = "'0.002935"
$ws.Range("E48").Value = "'-2.48%"
$ws.Range("E49").Value = "'-0.39%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.39%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.39%"
